# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets
# as produced by the latest scrape (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 2041
$wsExhibit.Range("F4").Value = 849
$wsExhibit.Range("F5").Value = 1126
$wsExhibit.Range("F6").Value = 350

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 2041
$wsAll.Range("F6").Value = 849
$wsAll.Range("F7").Value = 1126
$wsAll.Range("F8").Value = 350
